{"js": "// The document carries a stray \"_GoBack\" bookmark around the title\n// paragraph (left over from the last edit position). Since we are about to\n// type new content at the very end of the document, Word would normally\n// relocate that bookmark there instead -- so remove it from the title first\n// and re-add it after the newly typed text below.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst body = context.document.body;\nconst endRange = body.getRange(Word.RangeLocation.end);\n\n// Add the new \"P/S kazalnik\" bullet (continuing the same bulleted list --\n// numId 1 -- used by the other \"X kazalnik\" items above) together with its\n// explanatory paragraph, right after the existing last paragraph (\"EPS ...\").\nconst newContentOoxml =\n  '<?xml version=\"1.0\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"Odstavekseznama\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t>P/S kazalnik</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Ta kazalnik se uporablja za ocenjevanje mladih, potencialno hitro rasto\u010dih podjetij. Pove nam koliko so investitorji pripravljeni pla\u010dati na dolar prodajne vrednosti delnice.</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nendRange.insertOoxml(newContentOoxml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The original document carried a leftover \"_GoBack\" bookmark around the\n# title paragraph. When new content is typed at the end of the document,\n# Word moves that bookmark to mark the last edit location instead -- so we\n# drop it here and re-create it at the end, after the newly added text.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# Add the new \"P/S kazalnik\" bullet (same list/style as the other\n# \"X kazalnik\" bullets above it) plus its explanatory paragraph, appended\n# after the last paragraph in the body (the \"EPS\" one), before the sectPr.\n# NB: build the insertion point from a fresh, explicitly-collapsed Range\n# (not Paragraphs.Last.Range) -- InsertXML treats a non-collapsed range as\n# the span to overwrite, so reusing the paragraph's own Range would delete\n# the existing \"EPS\" paragraph text instead of appending after it.\n$endOfDoc = $d.Content.End\n$insertionRange = $d.Range($endOfDoc, $endOfDoc)\n\n$newParagraphsXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"Odstavekseznama\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>P/S kazalnik</w:t></w:r></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>Ta kazalnik se uporablja za ocenjevanje mladih, potencialno hitro rasto\u010dih podjetij. Pove nam koliko so investitorji pripravljeni pla\u010dati na dolar prodajne vrednosti delnice.</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n\n$insertionRange.InsertXML($newParagraphsXml)\n"}
